## feat: add 2022-Q1 data
##
## 1. Insert a new "2022-Q1" worksheet right before the "总计" (totals) sheet,
##    populated with the quarter's fund-holding detail rows (same layout as
##    the other quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/
##    仓位占比/持有市值(亿元)/仓位排名).
## 2. Prepend a "2022-Q1" summary row to the "总计" sheet (date / 持有数量(只) /
##    持有市值(亿元)), pushing the existing history rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" detail sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$totalsBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalsBefore)
$q1.Name = "2022-Q1"

# NOTE: inserting a sheet can shift what an existing worksheet variable
# resolves to (positional rebind), so re-look-up "总计" by name below
# instead of reusing $totalsBefore.

# Match the page-margin convention used by every other quarterly sheet.
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Pull over the header / row-1-data formatting (bold+border header style,
# centered index-column style) from the most recent quarter sheet so the
# new sheet renders consistently with its siblings.
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$prevQuarter.Range("A1:H2").Copy($q1.Range("A1"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
# Leading apostrophe forces these numeric-looking values to be stored as
# text, matching the source data (fund code / percentages kept as strings).
$q1.Range("B2").Value = "'004351"
$q1.Range("C2").Value = "汇丰晋信珠三角区域发展混合"
$q1.Range("D2").Value = "'0.51"
$q1.Range("E2").Value = "'93.92"
$q1.Range("F2").Value = "'5.65"
$q1.Range("G2").Value = "'0.0288"
$q1.Range("B2:G2").ClearFormats()
$q1.Range("H2").Value = 1

# The copy above also stamped A1 (blank in the source); clear it so the
# header row matches the rest of the workbook (no entry under the index
# column on row 1).
$q1.Range("A1").ClearContents()

# ---------------------------------------------------------------------
# 2) Prepend the 2022-Q1 totals row to "总计"
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()
$totals.Range("B2:D2").ClearFormats()
$totals.Range("A3").Copy($totals.Range("A2"))

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.03
